$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (Colaborador_id, Colaborador_nome, Departamento, Motivo_da_ausencia, Horas_de_ausencia, Data_da_ausencia, Salario)
$data = @(
    @(65100, "Alana Rezende", "P&D", "Problemas pessoais", 5, 45087, 2765.13),
    @(60989, "Srta. Ana Lívia Rodrigues", "Vendas", "Problemas pessoais", 4, 45079, 6709.26),
    @(60117, "Maria Sophia da Cruz", "Jurídico", "Doença", 3, 45098, 9167.530000000001),
    @(44273, "Dra. Maria Vitória Lopes", "Jurídico", "Problemas pessoais", 6, 45099, 9283.02),
    @(94046, "Sr. Danilo da Mota", "Marketing", "Viagem de negócios", 2, 45103, 12148.81),
    @(84679, "Cauê da Costa", "Engenharia", "Outros", 7, 45100, 7152.04),
    @(10291, "Fernando Barros", "Engenharia", "Doença", 2, 45096, 3200.58),
    @(31360, "Danilo Melo", "Marketing", "Doença", 2, 45093, 4037.76),
    @(60454, "Ana Júlia Martins", "Financeiro", "Consulta médica", 3, 45090, 4250.81),
    @(64853, "Isaac Dias", "P&D", "Viagem de negócios", 7, 45080, 9680.129999999999)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}
